# Edit: insert a new "2022-Q3" sheet (fund holdings for 科达利 002850)
# right after "总计", and update the "总计" summary sheet with the new
# quarter's row (shifting existing rows down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q3,
#    push existing rows down, and bump their running index in column A.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Re-apply the header/index style (bold, centered, bordered) that row 3
# (the old row 2, now shifted down) carries on column A, then set values.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(2,1).PasteSpecial(-4122)
$total.Cells.Item(2,1).Value = 0
$total.Range("B2:D2").ClearFormats()
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 38
$total.Cells.Item(2,4).Value = 4.82

# The rows below were shifted down by the insert but kept their old
# running index (column A) - renumber them (0-based position - 2).
for ($r = 3; $r -le 9; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" and fill it
#    with the fund-holding breakdown table.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}
# Match the header style used on every other quarter sheet (and on the
# "总计" sheet): bold, centered, bordered.
$total.Cells.Item(1,2).Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$rows = @(
"0|001606|农银汇理工业4.0灵活配置混合|40.20|86.85|3.43|1.3789|10",
"1|000336|农银研究精选混合|38.76|82.60|3.08|1.1938|7",
"2|481010|工银中小盘混合|15.67|87.16|2.76|0.4325|7",
"3|310308|申万菱信盛利精选混合|9.43|69.12|3.46|0.3263|8",
"4|004374|华泰保兴吉年丰混合A|5.25|93.40|5.23|0.2746|10",
"5|010149|浙商智选经济动能混合C|4.37|84.36|3.89|0.1700|8",
"6|005904|华泰保兴成长优选混合A|4.00|80.96|3.96|0.1584|6",
"7|013693|博道盛兴一年持有期混合|2.71|81.84|4.62|0.1252|4",
"8|012124|博道盛彦混合A|2.57|90.29|4.75|0.1221|5",
"9|004375|华泰保兴吉年丰混合C|2.26|93.40|5.23|0.1182|10",
"10|009847|圆信永丰研究精选混合A|1.14|89.51|4.71|0.0537|3",
"11|310368|申万菱信竞争优势混合A|1.05|92.76|4.54|0.0477|9",
"12|005933|新疆前海联合先进制造灵活配置混合A|0.95|92.14|4.88|0.0464|6",
"13|003493|申万菱信安鑫优选混合A|2.90|25.35|1.44|0.0418|5",
"14|009056|圆信永丰大湾区主题混合C|0.63|86.43|6.44|0.0406|3",
"15|006969|圆信永丰高端制造混合|0.91|87.79|3.60|0.0328|4",
"16|002210|创金合信量化多因子股票A|2.39|91.71|1.25|0.0299|6",
"17|003512|申万菱信安鑫优选混合C|1.98|25.35|1.44|0.0285|5",
"18|005009|申万菱信行业轮动股票A|0.62|92.87|4.38|0.0272|10",
"19|010148|浙商智选经济动能混合A|0.65|84.36|3.89|0.0253|8",
"20|009848|圆信永丰研究精选混合C|0.44|89.51|4.71|0.0207|3",
"21|009055|圆信永丰大湾区主题混合A|0.31|86.43|6.44|0.0200|3",
"22|012675|华融融泽6个月定开混合A|1.27|57.46|1.51|0.0192|2",
"23|003865|创金合信量化多因子股票C|0.75|91.71|1.25|0.0094|6",
"24|003647|创金合信中证1000指数增强C|0.68|90.54|1.25|0.0085|5",
"25|011761|平安鑫瑞混合A|0.59|23.80|1.42|0.0084|5",
"26|015707|安信新能源主题股票A|0.18|53.26|4.38|0.0079|4",
"27|011590|九泰天利量化股票C|0.50|83.77|1.47|0.0074|8",
"28|012125|博道盛彦混合C|0.13|90.29|4.75|0.0062|5",
"29|015708|安信新能源主题股票C|0.14|53.26|4.38|0.0061|4",
"30|005905|华泰保兴成长优选混合C|0.14|80.96|3.96|0.0055|6",
"31|015173|申万菱信竞争优势混合C|0.11|92.76|4.54|0.0050|9",
"32|005934|新疆前海联合先进制造灵活配置混合C|0.10|92.14|4.88|0.0049|6",
"33|003646|创金合信中证1000指数增强A|0.35|90.54|1.25|0.0044|5",
"34|011762|平安鑫瑞混合C|0.28|23.80|1.42|0.0040|5",
"35|012676|华融融泽6个月定开混合C|0.23|57.46|1.51|0.0035|2",
"36|015157|申万菱信行业轮动股票C|0.04|92.87|4.38|0.0018|10",
"37|011589|九泰天利量化股票A|0.07|83.77|1.47|0.0010|8"
)

# Column A (running index) and H (rank) are real numbers; B-G must stay
# text even though several of them look numeric (fund codes with leading
# zeros, percentages stored as strings, etc.) - mark that block as Text
# before writing so Excel does not "helpfully" convert it to a number.
$lastRow = 1 + $rows.Length
$q3.Range("B2:G$lastRow").NumberFormat = "@"

# Column A shares the same bold/centered/bordered style as the other
# quarter sheets' index column - apply it to the whole block up front
# (PasteSpecial only copies formatting, it does not touch values).
$total.Cells.Item(2,1).Copy()
$q3.Range("A2:A$lastRow").PasteSpecial(-4122)

$r = 2
foreach ($line in $rows) {
    $parts = $line.Split("|")
    $q3.Cells.Item($r, 1).Value = [double]$parts[0]
    $q3.Cells.Item($r, 2).Value = $parts[1]
    $q3.Cells.Item($r, 3).Value = $parts[2]
    $q3.Cells.Item($r, 4).Value = $parts[3]
    $q3.Cells.Item($r, 5).Value = $parts[4]
    $q3.Cells.Item($r, 6).Value = $parts[5]
    $q3.Cells.Item($r, 7).Value = $parts[6]
    $q3.Cells.Item($r, 8).Value = [double]$parts[7]
    $r++
}

Write-Output "done"
